$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RSTK-8172-New loc")
$ws.Activate()

# Delete entire row 17 (existing rows 18-22 shift up to become rows 17-21)
$ws.Rows.Item(17).Delete()

# Update selection on the sheet to match the post-delete state (row 17 selected, anchored at A17)
$ws.Range("A17:XFD17").Select()
